# Update "Total Submissions" figures (column B) for the rows that were
# corrected from placeholder 0s to the real submission counts, and roll
# the new per-person TOTAL rows forward to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2025-04-10
$ws.Range("B39").Value = 20
$ws.Range("B40").Value = 20
$ws.Range("B41").Value = 31

# 2025-04-11
$ws.Range("B43").Value = 16
$ws.Range("B44").Value = 15
$ws.Range("B45").Value = 42

# 2025-04-21
$ws.Range("B82").Value = 51
$ws.Range("B83").Value = 14
$ws.Range("B84").Value = 29
$ws.Range("B85").Value = 70

# 2025-04-22
$ws.Range("B86").Value = 24
$ws.Range("B87").Value = 16
$ws.Range("B88").Value = 12
$ws.Range("B89").Value = 40

# 2025-04-23
$ws.Range("B90").Value = 11
$ws.Range("B91").Value = 9
$ws.Range("B92").Value = 7
$ws.Range("B93").Value = 26

# TOTAL rows (per-person sums over the whole month)
$ws.Range("B122").Value = 110
$ws.Range("B123").Value = 253
$ws.Range("B124").Value = 425
$ws.Range("B125").Value = 143
